$d = $word.ActiveDocument

# 1) Replace the name "Fedrick James" -> "Pranay Raut"
$d.Content.Find.Execute("Fedrick James", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Pranay Raut", 2)

# 2) Replace the email-address paragraph built from three runs
#    "Fedric.james" + "23" + "@gmail.com" -> two runs "p" + "ranayraut12@gmail.com"
$d.Content.Find.Execute("Fedric.james23@gmail.com", $true, $false, $false, $false, $false,
                         $true, 1, $false, "pranayraut12@gmail.com", 2)

# 3) Merge "Jennifer" + ".james23@gmail.com" into a single run "Jennifer.james23@gmail.com"
$d.Content.Find.Execute("Jennifer.james23@gmail.com", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Jennifer.james23@gmail.com", 2)
